# Add the new "Kinect" purchase-request line (row 18) and its corresponding
# "Montant Par Personne" contribution entry (row 19), mirroring the existing
# "Moteur" row (row 17) formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New article label in A18, styled like the row above it (A17 / "Moteur").
$ws.Range("A18").Value = "Kinect"
$ws.Range("A17").Copy()
$ws.Range("A18").PasteSpecial(-4122)   # xlPasteFormats

# Price for the Kinect (130 Dt).
$ws.Range("E18").Value = 130

# Amin's contribution for the Kinect, tracked on the next row.
$ws.Range("K19").Value = 130

# Leave the selection on the newly added cell, like the saved workbook does.
$ws.Range("A18").Select()
